# Capital commitments workbook: add a new "Update Only" column (R) to the
# CapitalCommitment sheet, defaulted to "No" for every existing data row,
# and move the active selection/scroll position to reflect the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCommitment")

# Header for the new column.
$ws.Range("R1").Value = "Update Only"

# Default value for each of the 8 existing data rows.
$ws.Range("R2").Value = "No"
$ws.Range("R3").Value = "No"
$ws.Range("R4").Value = "No"
$ws.Range("R5").Value = "No"
$ws.Range("R6").Value = "No"
$ws.Range("R7").Value = "No"
$ws.Range("R8").Value = "No"
$ws.Range("R9").Value = "No"

# Match the author's final selection / scroll state (scrolled right so
# column D is the left-most visible column, with R9 the active cell).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("R9").Select()
